$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Unify all Password values (column B) to "UserUser123"
$ws.Range("B3").Value = "UserUser123"
$ws.Range("B4").Value = "UserUser123"
$ws.Range("B5").Value = "UserUser123"
$ws.Range("B6").Value = "UserUser123"
$ws.Range("B7").Value = "UserUser123"

# Update Result values (column E)
$ws.Range("E2").Value = "FAILED"
$ws.Range("E3").Value = "PASSES"
$ws.Range("E4").Value = "PASSES"
$ws.Range("E5").Value = "FAILED"
$ws.Range("E6").Value = "PASSES"
$ws.Range("E7").Value = "PASSES"

# Update the active selection to match the edited cell range
$ws.Range("E13").Select()
